# COVID Timeline.xlsx edit script
# Adds a new "TR" source row to the Sheet2 lookup table, then inserts a new
# timeline entry (row 75) on Sheet1 documenting the Central Valley Task Force /
# Unified Support Teams, shifting all later rows down by one.
# Finally restores the view-state (active tab / selections) to match the
# author's saved state: Sheet2 becomes the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Sheet2: add new source lookup row (A11/B11) -----------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A11").Value = "TR"
$ws2.Range("B11").Value = "Timeline provided in private communicaton from Raymundo, Trudy@CDPH"

# --- 2. Sheet1: insert new row 75 with the new timeline entry -------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(75).Insert() | Out-Null

$ws1.Range("A75").Value = 44040
$ws1.Range("B75").Value = "Unified Support Teams develpoed"
$ws1.Range("C75").Value = "declarations"
$ws1.Range("D75").Value = "TR"
$ws1.Range("E75").Value = 6
$ws1.Range("F75").Value = "Governor Gavin Newsom announced a call to action to slow the spread of COVID-19 in these hard-hit communities. On Tuesday, July 28, 2020, the California Governor’s Office of Emergency Services (Cal OES) and the California Health and Human Services Agency (CHHS) responded by establishing the Central Valley Task Force."

# --- 3. Restore view state --------------------------------------------------
$ws1.Range("A110").Select() | Out-Null

$ws2.Range("H20").Select() | Out-Null
$ws2.Activate() | Out-Null
